# Regenerate orders with updated distance/sizes.
# The experiment's Distance condition codes and the "S30" size code were
# renumbered (D64->D69, D80->D86, D51->D55, S30->S31). Apply the same
# substring substitution everywhere it appears: Condition, Filename_Left,
# Filename_Right, Distance and Size columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$firstCol = $usedRange.Column
$lastRow = $firstRow + $usedRange.Rows.Count - 1
$lastCol = $firstCol + $usedRange.Columns.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $orig = $cell.Value2

        if ($orig -eq $null) {
            continue
        }

        $updated = $orig
        $updated = $updated -replace "D64", "D69"
        $updated = $updated -replace "D80", "D86"
        $updated = $updated -replace "D51", "D55"
        $updated = $updated -replace "S30", "S31"

        if ($updated -ne $orig) {
            $cell.Value = $updated
        }
    }
}
